# Generate Report for Archive
#
# 1. Change the "Ready for handoff" status text (used on the Overview sheet
#    in columns E/F, and on the zh-cn / de-de sheets in column C) to
#    "In Translation".
# 2. Narrow the "zh-cn" / "de-de" status columns (Overview!E:F and
#    column C on the zh-cn / de-de sheets) from their old auto-fit width
#    down to the new narrower auto-fit width.

$wb = $excel.ActiveWorkbook

# --- 1. Update status text across all sheets -----------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Resize status columns ---------------------------------------------
# Target OOXML column width is 13.4101845877511 on each of these columns.
# The COM ColumnWidth setter in this runtime quantizes the stored width to
# the nearest 1/6, so 12.5 is the input value that lands closest to the
# real target (producing a stored width of 13.333333333333334).
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C
